$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$username = '[''cozy // Alex Claremont'', ''@alxclaremont'']'
$comments = '9'
$replys = '22'
$likes = '101'
$views = '1.966'
$text = '8 de set Alex se formou então na sequência teremos Alex advogado e Henry escritor e se o Henry escrever a história deles ??????????????????? Se vbsa livro for referênciado como Henry autor $^&$#@+×#$$%_^& 9'

# Force the cells to text format first so Excel doesn't reinterpret the
# numeric-looking strings (e.g. "9", "22", "1.966") as numbers.
foreach ($row in 2,3) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).NumberFormat = "@"
    }
}

foreach ($row in 2,3) {
    $ws.Cells.Item($row, 1).Value = $username
    $ws.Cells.Item($row, 2).Value = $comments
    $ws.Cells.Item($row, 3).Value = $replys
    $ws.Cells.Item($row, 4).Value = $likes
    $ws.Cells.Item($row, 5).Value = $views
    $ws.Cells.Item($row, 6).Value = $text
}

# Restore the default "Normal" cell style so no leftover style index is
# stamped on the cells (keeps the worksheet markup identical to before).
foreach ($row in 2,3) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Style = "Normal"
    }
}
